$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 381, pushing existing rows 381-397 down to 386-402
$ws.Range("A381:T385").EntireRow.Insert()

# Populate the 5 newly inserted rows (381-385) with new weekly data.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo,
# G Producto ID, H Producto, I Categoria ID, J Categoria, K Variedad,
# L Calidad, M Volumen, N Precio minimo, O Precio maximo,
# P Precio promedio ponderado, Q Unidad de comercializacion, R Origen,
# S Precio $/Kg, T Kg/unidad

$newRows = @(
    @{ Row = 381; D = 44509; K = "Clementina"; L = "Especial";                M = 300; N = 6000;  O = 6000;  P = 6000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota"; S = 600 },
    @{ Row = 382; D = 44509; K = "Clementina"; L = "Extra (doble especial)";  M = 410; N = 7000;  O = 7000;  P = 7000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota"; S = 700 },
    @{ Row = 383; D = 44509; K = "Clementina"; L = "Primera";                 M = 280; N = 5000;  O = 5000;  P = 5000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota"; S = 500 },
    @{ Row = 384; D = 44509; K = "Clementina"; L = "Segunda";                 M = 300; N = 4000;  O = 4000;  P = 4000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota"; S = 400 },
    @{ Row = 385; D = 44509; K = "Murcott";    L = "Tercera";                 M = 350; N = 3000;  O = 3000;  P = 3000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de San Felipe de Aconcagua"; S = 300 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102004
    $ws.Cells.Item($row, 10).Value = "Mandarina"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}
